$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.660.51"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "1.845.02"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5285"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3162"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06803"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07769"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "1.854.69"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.015"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9993"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007921"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "26.688.21"
$ws.Range("E20").Value = "  +0.93%  "

$ws.Range("D21").Value = "2.078.10"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.985"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.329"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.223"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.677"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08724"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.091"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04898"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7312"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.95%  "

$ws.Range("E35").Value = "  +1.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.852"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.284"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01731"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  -1.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.958"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.704"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4184"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.079"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1240"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.04%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05811"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8912"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "
